$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.673.12"
$ws.Range("E2").Value = "  +2.21%  "

$ws.Range("D3").Value = "1.870.36"
$ws.Range("E3").Value = "  +1.13%  "

$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").Value = "'312.86"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").Value = "'1.012"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("D7").Value = "'0.4794"
$ws.Range("E7").Value = "  +0.64%  "

$ws.Range("D8").Value = "'0.3815"
$ws.Range("E8").Value = "  +3.51%  "

$ws.Range("D9").Value = "'0.07369"
$ws.Range("E9").Value = "  +1.77%  "

$ws.Range("D10").Value = "'0.9366"
$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").Value = "'20.88"
$ws.Range("E11").Value = "  +4.89%  "

$ws.Range("D12").Value = "'0.07812"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").Value = "1.883.74"
$ws.Range("E13").Value = "  +1.27%  "

$ws.Range("D14").Value = "'5.457"
$ws.Range("E14").Value = "  +1.09%  "

$ws.Range("D15").Value = "'6.580"
$ws.Range("E15").Value = "  +1.46%  "

$ws.Range("D16").Value = "'90.57"
$ws.Range("E16").Value = "  +1.63%  "

$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").Value = "'0.000008855"
$ws.Range("E18").Value = "  +2.08%  "

$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").Value = "27.608.55"
$ws.Range("E20").Value = "  +1.78%  "

$ws.Range("D21").Value = "'14.73"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("D22").Value = "'5.111"
$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("D23").Value = "'10.76"
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("D24").Value = "'1.938"
$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("D25").Value = "'156.51"
$ws.Range("E25").Value = "  +2.30%  "

$ws.Range("D26").Value = "'18.53"
$ws.Range("E26").Value = "  +0.90%  "

$ws.Range("D27").Value = "'2.035"
$ws.Range("E27").Value = "  +2.35%  "

$ws.Range("D28").Value = "'115.75"
$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("D29").Value = "'4.960"
$ws.Range("E29").Value = "  +0.74%  "

$ws.Range("D30").Value = "'0.08895"
$ws.Range("E30").Value = "  +0.34%  "

$ws.Range("D31").Value = "'3.331"
$ws.Range("E31").Value = "  +0.33%  "

$ws.Range("D32").Value = "'1.215"
$ws.Range("E32").Value = "  +3.02%  "

$ws.Range("D33").Value = "'0.7629"
$ws.Range("E33").Value = "  +3.43%  "

$ws.Range("D34").Value = "'4.611"
$ws.Range("E34").Value = "  +2.04%  "

$ws.Range("E35").Value = "  +1.07%  "

$ws.Range("D36").Value = "'1.134"
$ws.Range("E36").Value = "  +1.71%  "

$ws.Range("D37").Value = "'0.02039"
$ws.Range("E37").Value = "  +3.13%  "

$ws.Range("D38").Value = "'0.5687"

$ws.Range("D39").Value = "'0.05393"
$ws.Range("E39").Value = "  +2.50%  "

$ws.Range("D40").Value = "'2.983"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").Value = "'7.061"
$ws.Range("E41").Value = "  +0.23%  "

$ws.Range("D42").Value = "'8.570"
$ws.Range("E42").Value = "  +3.35%  "

$ws.Range("D43").Value = "'0.1529"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'10.76"
$ws.Range("E44").Value = "  +1.50%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.4906"
$ws.Range("E45").Value = "  +3.43%  "

$ws.Range("D46").Value = "'105.17"
$ws.Range("E46").Value = "  +3.08%  "

$ws.Range("D47").Value = "'1.014"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").Value = "'1.667"
$ws.Range("E48").Value = "  +3.26%  "

$ws.Range("D49").Value = "'67.60"
$ws.Range("E49").Value = "  +2.34%  "

$ws.Range("D50").Value = "'0.06106"
$ws.Range("E50").Value = "  +0.75%  "

$ws.Range("D51").Value = "'0.9119"
$ws.Range("E51").Value = "  +2.04%  "
